$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (e.g. "13.00", "6.10") stay text,
# matching the source data's inlineStr representation instead of being
# auto-coerced to numbers by Excel (which would also strip trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.558.90'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '3.744.48'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '602.02'
$ws.Range("E5").Value = '  +5.44%  '
$ws.Range("D6").Value = '186.38'
$ws.Range("E6").Value = '  +15.98%  '
$ws.Range("D7").Value = '3.737.45'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.636'
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '0.727'
$ws.Range("E10").Value = '  +0.96%  '
$ws.Range("D11").Value = '0.164'
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D12").Value = '57.24'
$ws.Range("E12").Value = '  +12.50%  '
$ws.Range("D13").Value = '0.0000296'
$ws.Range("E13").Value = '  -3.65%  '
$ws.Range("D14").Value = '10.86'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '4.350.99'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '3.761.98'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '19.58'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.126'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '13.00'
$ws.Range("E19").Value = '  -2.30%  '
$ws.Range("D20").Value = '1.13'
$ws.Range("E20").Value = '  -3.05%  '
$ws.Range("D21").Value = '69.353.60'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").Value = '415.05'
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("D24").Value = '89.69'
$ws.Range("D25").Value = '3.07'
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").Value = '12.96'
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("D27").Value = '11.07'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("D28").Value = '3.99'
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("D29").Value = '6.10'
$ws.Range("E29").Value = '  +3.02%  '
$ws.Range("D30").Value = '9.58'
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("D31").Value = '33.13'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").Value = '7.39'
$ws.Range("E32").Value = '  -4.62%  '
$ws.Range("D33").Value = '12.57'
$ws.Range("E33").Value = '  -3.35%  '
$ws.Range("D34").Value = '0.119'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").Value = '44.29'
$ws.Range("E35").Value = '  -3.51%  '
$ws.Range("D36").Value = '615.79'
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = '65.38'
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("D38").Value = '0.0₃0909'
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").Value = '0.408'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = '0.138'
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("D43").Value = '3.08'
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  +1.85%  '
$ws.Range("D45").Value = '3.00'
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("D46").Value = '0.0446'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("D47").Value = '9.36'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.136'
$ws.Range("E48").Value = '  -2.10%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.795.66'
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").Value = '3.22'
$ws.Range("E50").Value = '  +1.55%  '
$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = '2.73'
$ws.Range("E51").Value = '  -0.83%  '
